# Insert a new row at position 499 (pushing the existing rows 499-524
# down to 500-525) and populate it with the new "A_SERIES_POSTSEC" entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 499; everything from 499 downward shifts by one.
$ws.Rows("499:499").Insert()

# Copy the formatting (style s="4") of the data row immediately below the
# newly inserted blank row so the new row matches the look of its
# neighbours, then set the actual values.
$ws.Range("A500:D500").Copy()
$ws.Range("A499:D499").PasteSpecial(-4122)

$ws.Range("A499").Value = "A_SERIES_POSTSEC"
$ws.Range("B499").Value = "K_SERIES"
$ws.Range("C499").Value = "Postsekundarer nicht-tertiärer Abschluss"
$ws.Range("D499").Value = "With post-secondary non-tertiary qualification"
